# Scheduled-runner update: refresh Market Board price snapshots + recompute
# leve profit columns (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) for the affected leves across the ALC, BSM, CRP, CUL,
# GSM and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 485.33334
$ws.Range("I18").Value = 328
$ws.Range("J18").Value = 800
$ws.Range("K18").Value = 328
$ws.Range("L18").Value = 800
$ws.Range("M18").Value = -44
$ws.Range("N18").Value = -1368
$ws.Range("H45").Value = 2875
$ws.Range("I45").Value = 2875
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 8625
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -8433
$ws.Range("N45").ClearContents()
$ws.Range("H64").Value = 3808.1738
$ws.Range("J64").Value = 4487.875
$ws.Range("L64").Value = 4487.875
$ws.Range("N64").Value = -4983.875
$ws.Range("H67").Value = 3808.1738
$ws.Range("J67").Value = 4487.875
$ws.Range("L67").Value = 4487.875
$ws.Range("N67").Value = -6203.875
$ws.Range("H76").Value = 3857.8684
$ws.Range("I76").Value = 3681.9644
$ws.Range("J76").Value = 4350.4
$ws.Range("K76").Value = 3681.9644
$ws.Range("L76").Value = 4350.4
$ws.Range("M76").Value = -3366.9644
$ws.Range("N76").Value = -4980.4
$ws.Range("H79").Value = 3857.8684
$ws.Range("I79").Value = 3681.9644
$ws.Range("J79").Value = 4350.4
$ws.Range("K79").Value = 3681.9644
$ws.Range("L79").Value = 4350.4
$ws.Range("M79").Value = -2589.9644
$ws.Range("N79").Value = -6534.4
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H88").Value = 1525.375
$ws.Range("I88").Value = 1233.8334
$ws.Range("J88").Value = 2400
$ws.Range("K88").Value = 1233.8334
$ws.Range("L88").Value = 2400
$ws.Range("M88").Value = -827.8334
$ws.Range("N88").Value = -3212
$ws.Range("H91").Value = 1525.375
$ws.Range("I91").Value = 1233.8334
$ws.Range("J91").Value = 2400
$ws.Range("K91").Value = 1233.8334
$ws.Range("L91").Value = 2400
$ws.Range("M91").Value = 170.1666
$ws.Range("N91").Value = -5208
$ws.Range("H97").Value = 650
$ws.Range("I97").Value = 200
$ws.Range("J97").Value = 1100
$ws.Range("K97").Value = 600
$ws.Range("L97").Value = 3300
$ws.Range("M97").Value = -104
$ws.Range("N97").Value = -4292
$ws.Range("H100").Value = 1690.25
$ws.Range("I100").Value = 1687.5
$ws.Range("J100").Value = 1698.5
$ws.Range("K100").Value = 1687.5
$ws.Range("L100").Value = 1698.5
$ws.Range("M100").Value = -1146.5
$ws.Range("N100").Value = -2780.5
$ws.Range("H106").Value = 3550.4443
$ws.Range("J106").Value = 4640
$ws.Range("L106").Value = 4640
$ws.Range("N106").Value = -5902
$ws.Range("H109").Value = 58796
$ws.Range("J109").Value = 58796
$ws.Range("L109").Value = 58796
$ws.Range("N109").Value = -61570
$ws.Range("H112").Value = 5008
$ws.Range("J112").Value = 1670.8572
$ws.Range("L112").Value = 5012.571599999999
$ws.Range("N112").Value = -7228.571599999999
$ws.Range("H115").Value = 2457
$ws.Range("I115").Value = 2457
$ws.Range("K115").Value = 7371
$ws.Range("M115").Value = -5804
$ws.Range("H137").Value = 2910.68
$ws.Range("I137").Value = 3288.6316
$ws.Range("J137").Value = 1713.8334
$ws.Range("K137").Value = 9865.8948
$ws.Range("L137").Value = 5141.5002
$ws.Range("M137").Value = -7315.8948
$ws.Range("N137").Value = -10241.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 671.9048
$ws.Range("I94").Value = 703.5714
$ws.Range("J94").Value = 608.5714
$ws.Range("K94").Value = 703.5714
$ws.Range("L94").Value = 608.5714
$ws.Range("M94").Value = -252.5714
$ws.Range("N94").Value = -1510.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 60007.5
$ws.Range("J23").Value = 60007.5
$ws.Range("L23").Value = 60007.5
$ws.Range("N23").Value = -60487.5
$ws.Range("H27").Value = 60007.5
$ws.Range("J27").Value = 60007.5
$ws.Range("L27").Value = 60007.5
$ws.Range("N27").Value = -60391.5
$ws.Range("H31").Value = 8445.549999999999
$ws.Range("I31").Value = 7560.7334
$ws.Range("K31").Value = 7560.7334
$ws.Range("M31").Value = -7265.7334
$ws.Range("H34").Value = 8445.549999999999
$ws.Range("I34").Value = 7560.7334
$ws.Range("K34").Value = 7560.7334
$ws.Range("M34").Value = -7358.7334
$ws.Range("H39").Value = 7525.5
$ws.Range("I39").Value = 7525.5
$ws.Range("K39").Value = 7525.5
$ws.Range("M39").Value = -7134.5
$ws.Range("H49").Value = 7525.5
$ws.Range("I49").Value = 7525.5
$ws.Range("K49").Value = 7525.5
$ws.Range("M49").Value = -7343.5
$ws.Range("H62").Value = 3710.75
$ws.Range("I62").Value = 3572.1428
$ws.Range("J62").Value = 3904.8
$ws.Range("K62").Value = 3572.1428
$ws.Range("L62").Value = 3904.8
$ws.Range("M62").Value = -2948.1428
$ws.Range("N62").Value = -5152.8
$ws.Range("H65").Value = 3710.75
$ws.Range("I65").Value = 3572.1428
$ws.Range("J65").Value = 3904.8
$ws.Range("K65").Value = 17860.714
$ws.Range("L65").Value = 19524
$ws.Range("M65").Value = -14740.714
$ws.Range("N65").Value = -25764
$ws.Range("H109").Value = 39500
$ws.Range("J109").Value = 39500
$ws.Range("L109").Value = 39500
$ws.Range("N109").Value = -41580

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 2000
$ws.Range("J35").Value = 2000
$ws.Range("L35").Value = 6000
$ws.Range("N35").Value = -6576
$ws.Range("H122").Value = 1302.1538
$ws.Range("J122").Value = 1620.4
$ws.Range("L122").Value = 14583.6
$ws.Range("N122").Value = -19483.6
$ws.Range("H133").Value = 3400.9
$ws.Range("I133").Value = 2966.5
$ws.Range("J133").Value = 3835.3
$ws.Range("K133").Value = 8899.5
$ws.Range("L133").Value = 11505.9
$ws.Range("M133").Value = -3839.5
$ws.Range("N133").Value = -21625.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 20166.666
$ws.Range("I43").Value = 500
$ws.Range("J43").Value = 30000
$ws.Range("K43").Value = 500
$ws.Range("L43").Value = 30000
$ws.Range("M43").Value = -349
$ws.Range("N43").Value = -30302

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 57266.332
$ws.Range("J27").Value = 57266.332
$ws.Range("L27").Value = 57266.332
$ws.Range("N27").Value = -57404.332
$ws.Range("H96").Value = 4289.1113
$ws.Range("I96").Value = 2149.5
$ws.Range("J96").Value = 6000.8
$ws.Range("K96").Value = 2149.5
$ws.Range("L96").Value = 6000.8
$ws.Range("M96").Value = -776.5
$ws.Range("N96").Value = -8746.799999999999
$ws.Range("H109").Value = 70400
$ws.Range("J109").Value = 70400
$ws.Range("L109").Value = 70400
$ws.Range("N109").Value = -73174
$ws.Range("H115").Value = 74900
$ws.Range("J115").Value = 74900
$ws.Range("L115").Value = 74900
$ws.Range("N115").Value = -78034

